# Applies the edits described by the commit:
#   1. Bump the auto-generated "last displayed" date field from 4/4/11 to
#      4/19/11 everywhere it is baked into the deck (slide master + every
#      slide layout).
#   2. On the slide's "TextBox 32" shape (the ArrayList size annotation):
#        - narrow the shape slightly (cx 1745878 -> 1730449 EMU)
#        - split "X100,000 = 7.629MB" into three runs ("x", "100,000 ",
#          "= 7.629MB"), lower-casing the leading X, matching the sibling
#          "Product" textbox's run layout.

$p = $ppt.ActivePresentation

function Update-DateField($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.Name -like "Date Placeholder*") {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "4/4/11") {
                $tr.Text = "4/19/11"
            }
        }
    }
}

# --- 1. Update the date placeholder on the master and on every layout ---
$master = $p.SlideMaster
Update-DateField $master.Shapes

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $lyt = $master.CustomLayouts.Item($i)
    Update-DateField $lyt.Shapes
}

# --- 2. Update the "ArrayList" annotation textbox on slide 1 ---
$s = $p.Slides.Item(1)

for ($k = 1; $k -le $s.Shapes.Count; $k++) {
    $candidate = $s.Shapes.Item($k)
    if ($candidate.Name -eq "TextBox 32") {
        $shp = $candidate
    }
}

# Narrow the shape to its new width (EMU -> points, 12700 EMU per point).
$shp.Width = 1730449 / 12700

$tr = $shp.TextFrame.TextRange
$fullText = $tr.Text
$idx = $fullText.IndexOf("X100,000 = 7.629MB")

if ($idx -ge 0) {
    $start = $idx + 1

    # "X" -> "x"
    $tr.Characters($start, 1).Text = "x"

    # "100,000 " stays the same text, but becomes its own run
    $tr.Characters($start + 1, 8).Text = "100,000 "

    # "= 7.629MB" stays the same text, but becomes its own run
    $tr.Characters($start + 9, 9).Text = "= 7.629MB"
}
